$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing row (31) currently has the "last-row" date style (YYYY-MM-DD).
# Once we append a new row, row 31 becomes a regular row and should pick up the
# regular date+time style used by the rest of column A (copy it from A30, which
# already has that "regular" style).
$ws.Range("A30").Copy()
$ws.Range("A31").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Append the new day's data as row 32, and give A32 the distinctive "last-row"
# date style previously used by A31 (copy formatting from old A31 style source).
$ws.Range("A32").Value = 45772
$ws.Range("A32").NumberFormat = "YYYY-MM-DD"
$ws.Range("B32").Value = 130
$ws.Range("C32").Value = 129
$ws.Range("D32").Value = 129
